# Update epexspot_prices.xlsx with the latest day of data:
#   - "Prix Spot": a new date column (31-aug) with 24 hourly prices.
#   - "Gaz" / "CO2": a new row (2025-08-29) with the day's closing price.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append column CA (31-aug) right after BZ (30-aug).
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (BZ1) onto the new header
# cell so CA1 matches the existing header styling (bold, bordered, centered)
# instead of picking up the default style.
$wsSpot.Range("BZ1").Copy() | Out-Null
$wsSpot.Range("CA1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$wsSpot.Range("CA1").Value2 = "31-aug"

$caValues = @(
    15.86,
    12.93,
    12.33,
    7.61,
    5.95,
    8.119999999999999,
    9.18,
    16.24,
    22.15,
    3.6,
    0,
    -0.01,
    -0.08,
    -0.95,
    -0.6,
    -0.01,
    4.22,
    4.96,
    17.53,
    47.15,
    59.55,
    73.88,
    71.40000000000001,
    62
)

for ($i = 0; $i -lt $caValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 79).Value2 = $caValues[$i]   # column 79 = CA
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append the latest day's closing price as row 76.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date column to stay plain text (matching every existing row,
# which stores dates as literal strings, not real Excel dates) — otherwise
# Excel's automatic type inference on a "YYYY-MM-DD"-looking string turns it
# into a date serial number. Apply a text format before assigning, then
# clear formatting again so the cell doesn't end up with a lingering
# "@"/date style that the other date cells in the column don't have.
$wsGaz.Range("A76").NumberFormat = "@"
$wsGaz.Range("A76").Value2 = "2025-08-29"
$wsGaz.Range("A76").ClearFormats() | Out-Null

$wsGaz.Range("B76").Value2 = 30.375

# ---------------------------------------------------------------------------
# Sheet "CO2": append the latest day's closing price as row 76.
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A76").NumberFormat = "@"
$wsCO2.Range("A76").Value2 = "2025-08-29"
$wsCO2.Range("A76").ClearFormats() | Out-Null

$wsCO2.Range("B76").Value2 = 71.09999999999999
